$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" '69.494.63'
Set-TextValue "E2" '  -0.59%  '
Set-TextValue "D3" '2.489.66'
Set-TextValue "E3" '  -1.11%  '
Set-TextValue "E4" '  -0.10%  '
Set-TextValue "D5" '568.23'
Set-TextValue "E5" '  -1.37%  '
Set-TextValue "D6" '164.30'
Set-TextValue "E6" '  -1.45%  '
Set-TextValue "E7" '  -0.10%  '
Set-TextValue "D8" '0.512'
Set-TextValue "E8" '  -1.22%  '
Set-TextValue "D9" '2.488.60'
Set-TextValue "E9" '  -0.97%  '
Set-TextValue "E10" '  -2.16%  '
Set-TextValue "E11" '  -0.43%  '
Set-TextValue "D12" '0.354'
Set-TextValue "E12" '  -0.30%  '
Set-TextValue "E13" '  -0.77%  '
Set-TextValue "D14" '2.946.80'
Set-TextValue "E14" '  -1.55%  '
Set-TextValue "D15" '69.362.19'
Set-TextValue "E15" '  -0.75%  '
Set-TextValue "E16" '  -0.86%  '
Set-TextValue "D17" '24.22'
Set-TextValue "E17" '  -3.15%  '
Set-TextValue "D18" '2.499.04'
Set-TextValue "E18" '  -0.53%  '
Set-TextValue "D19" '11.15'
Set-TextValue "E19" '  -2.52%  '
Set-TextValue "D20" '7.37'
Set-TextValue "E20" '  -5.38%  '
Set-TextValue "D21" '346.63'
Set-TextValue "E21" '  -1.14%  '
Set-TextValue "E22" '  -1.38%  '
Set-TextValue "E23" '  -4.98%  '
Set-TextValue "E24" '  -0.08%  '
Set-TextValue "D25" '69.57'
Set-TextValue "E25" '  -1.24%  '
Set-TextValue "D26" '3.90'
Set-TextValue "E26" '  -2.74%  '
Set-TextValue "D27" '2.619.50'
Set-TextValue "E27" '  -1.04%  '
Set-TextValue "E28" '  -3.50%  '
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  +0.02%  '
Set-TextValue "D30" '0.0₃0869'
Set-TextValue "E31" '  -4.08%  '
Set-TextValue "E32" '  -5.47%  '
Set-TextValue "D33" '436.77'
Set-TextValue "E33" '  -6.21%  '
Set-TextValue "D34" '0.999'
Set-TextValue "E34" '  -0.11%  '
Set-TextValue "E35" '  -2.03%  '
Set-TextValue "E36" '  -0.27%  '
Set-TextValue "D37" '19.07'
Set-TextValue "E38" '  -3.88%  '
Set-TextValue "D39" '18.15'
Set-TextValue "E39" '  -2.42%  '
Set-TextValue "E40" '  +0.01%  '
Set-TextValue "E41" '  -1.76%  '
Set-TextValue "E42" '  -4.21%  '
Set-TextValue "E43" '  -2.02%  '
Set-TextValue "E44" '  -6.07%  '
Set-TextValue "E45" '  -6.96%  '
Set-TextValue "D46" '138.80'
Set-TextValue "E46" '  -2.72%  '
Set-TextValue "E47" '  -2.04%  '
Set-TextValue "D48" '0.510'
Set-TextValue "E48" '  -3.61%  '
Set-TextValue "E49" '  -1.42%  '
Set-TextValue "E50" '  -0.72%  '
Set-TextValue "D51" '1.93'
Set-TextValue "E51" '  +19.54%  '
